$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (1-based, including header row 1) that correspond to recipes which
# contain allergy-triggering ingredients (brinjal/eggplant, mushrooms, okra,
# eggs) and must be removed entirely from the sheet.
$rowsToDelete = @(41, 39, 13, 5)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# For every remaining recipe row, the "Preparation method" column (H) gets a
# leading ", " and loses its trailing padding whitespace.
$usedRows = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $usedRows; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = ", " + $val.TrimEnd()
        $cell.Value = $newVal
    }
}
